$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B, shifting the existing "title" column (and
# everything to its right) one column to the right.
$ws.Columns.Item(2).Insert()

# Populate the newly inserted column with the "dateCreated" header and
# its two numeric values.
$ws.Range("B1").Value = "dateCreated"
$ws.Range("B2").Value = 1234
$ws.Range("B3").Value = 2222

# Match the author's final selection.
[void]$ws.Range("B3").Select()
